# Update TPM-derived ligand/receptor expression statistics (new TPM run)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 3.045497666666666
$ws.Range("H2").Value2 = 9.136493
$ws.Range("I2").Value2 = 0.06184575966423571
$ws.Range("J2").Value2 = 0.06184575966423572
$ws.Range("O2").Value2 = 0.8416031693647025
$ws.Range("P2").Value2 = 0.8416031693647025
$ws.Range("Q2").Value2 = 4.809998104779999
$ws.Range("R2").Value2 = 43.28998294301999
$ws.Range("S2").Value2 = 0.05204958734518845
$ws.Range("T2").Value2 = 0.05204958734518846

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 3.045497666666666
$ws.Range("H3").Value2 = 9.136493
$ws.Range("I3").Value2 = 0.06184575966423571
$ws.Range("J3").Value2 = 0.06184575966423572
$ws.Range("M3").Value2 = 0.2972526666666667
$ws.Range("N3").Value2 = 0.891758
$ws.Range("O3").Value2 = 0.1583968306352975
$ws.Range("P3").Value2 = 0.1583968306352975
$ws.Range("Q3").Value2 = 0.9052823027437777
$ws.Range("R3").Value2 = 8.147540724694
$ws.Range("S3").Value2 = 0.009796172319047256
$ws.Range("T3").Value2 = 0.009796172319047258

# Row 4
$ws.Range("I4").Value2 = 0.6352626115862781
$ws.Range("J4").Value2 = 0.6352626115862781
$ws.Range("O4").Value2 = 0.8416031693647025
$ws.Range("P4").Value2 = 0.8416031693647025
$ws.Range("S4").Value2 = 0.5346390272899095
$ws.Range("T4").Value2 = 0.5346390272899095

# Row 5
$ws.Range("I5").Value2 = 0.6352626115862781
$ws.Range("J5").Value2 = 0.6352626115862781
$ws.Range("M5").Value2 = 0.2972526666666667
$ws.Range("N5").Value2 = 0.891758
$ws.Range("O5").Value2 = 0.1583968306352975
$ws.Range("P5").Value2 = 0.1583968306352975
$ws.Range("Q5").Value2 = 9.29881050836889
$ws.Range("R5").Value2 = 83.68929457532002
$ws.Range("S5").Value2 = 0.1006235842963685
$ws.Range("T5").Value2 = 0.1006235842963685

# Row 6
$ws.Range("G6").Value2 = 14.91542433333333
$ws.Range("H6").Value2 = 44.746273
$ws.Range("I6").Value2 = 0.3028916287494862
$ws.Range("J6").Value2 = 0.3028916287494862
$ws.Range("O6").Value2 = 0.8416031693647025
$ws.Range("P6").Value2 = 0.8416031693647025
$ws.Range("Q6").Value2 = 23.55712288358
$ws.Range("R6").Value2 = 212.01410595222
$ws.Range("S6").Value2 = 0.2549145547296044
$ws.Range("T6").Value2 = 0.2549145547296044

# Row 7
$ws.Range("G7").Value2 = 14.91542433333333
$ws.Range("H7").Value2 = 44.746273
$ws.Range("I7").Value2 = 0.3028916287494862
$ws.Range("J7").Value2 = 0.3028916287494862
$ws.Range("M7").Value2 = 0.2972526666666667
$ws.Range("N7").Value2 = 0.891758
$ws.Range("O7").Value2 = 0.1583968306352975
$ws.Range("P7").Value2 = 0.1583968306352975
$ws.Range("Q7").Value2 = 4.433649657548222
$ws.Range("R7").Value2 = 39.90284691793401
$ws.Range("S7").Value2 = 0.04797707401988178
$ws.Range("T7").Value2 = 0.04797707401988178

